$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G10/H10 currently hold numeric 0; replace with the text "zeros".
# Copy the number-format/font styling already used by F10 (style index 1)
# onto G10:H10 first, then overwrite the values, so the cells end up
# styled + typed exactly like the author's edit (shared string "zeros").
$ws.Range("F10").Copy() | Out-Null
$ws.Range("G10:H10").PasteSpecial(-4122) | Out-Null
$ws.Range("G10").Value = "zeros"
$ws.Range("H10").Value = "zeros"

# G11/H11/I11/J11 currently hold numeric 0; replace with text "zeros",
# using F11's styling (also style index 1) the same way.
$ws.Range("F11").Copy() | Out-Null
$ws.Range("G11:J11").PasteSpecial(-4122) | Out-Null
$ws.Range("G11").Value = "zeros"
$ws.Range("H11").Value = "zeros"
$ws.Range("I11").Value = "zeros"
$ws.Range("J11").Value = "zeros"

$excel.CutCopyMode = 0

# Move the visible selection to H17, matching the saved cursor position.
$ws.Range("H17").Select() | Out-Null
